$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.026.82'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '2.548.61'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range("D5").Value = "'583.24"
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range("D6").Value = "'146.73"
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range("D10").Value = "'5.55"
$ws.Range('E10').Value = '  -3.48%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '62.925.94'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').Value = '2.543.49'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range("D19").Value = "'338.09"
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range("D21").Value = "'6.75"
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range("D23").Value = "'65.63"
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range("D29").Value = "'8.34"
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range("D30").Value = "'7.68"
$ws.Range('E30').Value = '  +6.31%  '
$ws.Range("D31").Value = "'1.97"
$ws.Range('E31').Value = '  +4.65%  '
$ws.Range('D32').Value = '0.0₃0815'
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range("D33").Value = "'177.97"
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range("D34").Value = "'1.54"
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range("D35").Value = "'417.44"
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range("D37").Value = "'19.10"
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('E40').Value = '  -2.26%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range("D42").Value = "'39.77"
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range("D43").Value = "'151.07"
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range("D45").Value = "'20.76"
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range("D47").Value = "'0.603"
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range("D50").Value = "'18.31"
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('E51').Value = '  -6.35%  '
